# Addition of filtering option + fixing orientation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "revenue" column header
$ws.Range("E1").Value = "revenue"

# Revenue values for each data row (rows 2-23)
$revenue = @{
    2  = 100
    3  = 125
    4  = 150
    5  = 900
    6  = 700
    7  = 1200
    8  = 111
    9  = 90
    10 = 400
    11 = 600
    12 = 800
    13 = 888
    14 = 123
    15 = 145
    16 = 1556
    17 = 1720
    18 = 677
    19 = 870
    20 = 875
    21 = 345
    22 = 90
    23 = 1200
}

foreach ($row in $revenue.Keys) {
    $ws.Cells.Item($row, 5).Value = $revenue[$row]
}

# Fix the orientation / selection of the active cell
$ws.Range("C3").Select() | Out-Null
